# Update cryptocurrency price/volume data for Wed Mar 15 16:20:15 UTC 2023 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '24.420.97'
$ws.Range('E2').Value = '  -6.03%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '1.650.42'
$ws.Range('E3').Value = '  -6.73%  '

# Row 4: TetherUSD
$ws.Range('D4').Value = '''1.013'
$ws.Range('E4').Value = '  +1.40%  '

# Row 5: USDC
$ws.Range('D5').Value = '''1.007'
$ws.Range('E5').Value = '  +1.15%  '

# Row 6: BNB
$ws.Range('D6').Value = '''307.33'
$ws.Range('E6').Value = '  -3.08%  '

# Row 7: XRP
$ws.Range('D7').Value = '''0.3637'
$ws.Range('E7').Value = '  -5.34%  '

# Row 8: OKB
$ws.Range('D8').Value = '''46.92'
$ws.Range('E8').Value = '  -7.94%  '

# Row 9: Cardano
$ws.Range('D9').Value = '''0.3264'
$ws.Range('E9').Value = '  -10.35%  '

# Row 10: Polygon
$ws.Range('D10').Value = '''1.110'
$ws.Range('E10').Value = '  -10.37%  '

# Row 11: Dogecoin
$ws.Range('D11').Value = '''0.06949'
$ws.Range('E11').Value = '  -9.38%  '

# Row 12: BinanceUSD
$ws.Range('D12').Value = '''1.012'
$ws.Range('E12').Value = '  +1.72%  '

# Row 13: Polkadot
$ws.Range('D13').Value = '''5.977'
$ws.Range('E13').Value = '  -8.25%  '

# Row 14: Solana
$ws.Range('D14').Value = '''19.38'
$ws.Range('E14').Value = '  -11.59%  '

# Row 15: Chainlink
$ws.Range('D15').Value = '''6.608'
$ws.Range('E15').Value = '  -7.19%  '

# Row 16: WrappedEther
$ws.Range('D16').Value = '1.645.94'
$ws.Range('E16').Value = '  -6.99%  '

# Row 17: ShibaInu
$ws.Range('D17').Value = '''0.00001047'
$ws.Range('E17').Value = '  -10.05%  '

# Row 18: TRON
$ws.Range('D18').Value = '''0.06584'
$ws.Range('E18').Value = '  -3.90%  '

# Row 19: Dai
$ws.Range('D19').Value = '''1.006'
$ws.Range('E19').Value = '  +1.18%  '

# Row 20: Litecoin
$ws.Range('D20').Value = '''78.07'
$ws.Range('E20').Value = '  -10.85%  '

# Row 21: Avalanche
$ws.Range('D21').Value = '''15.92'
$ws.Range('E21').Value = '  -10.32%  '

# Row 22: Uniswap
$ws.Range('D22').Value = '''5.959'

# Row 23: Cosmos
$ws.Range('D23').Value = '''12.01'
$ws.Range('E23').Value = '  -6.47%  '

# Row 24: WrappedBTC
$ws.Range('D24').Value = '24.485.01'
$ws.Range('E24').Value = '  -5.56%  '

# Row 25: Toncoin
$ws.Range('D25').Value = '''2.424'

# Row 26: LidoDAOToken
$ws.Range('D26').Value = '''2.426'
$ws.Range('E26').Value = '  -19.03%  '

# Row 27: Monero
$ws.Range('D27').Value = '''147.00'
$ws.Range('E27').Value = '  -5.83%  '

# Row 28: EthereumClassic
$ws.Range('D28').Value = '''18.98'
$ws.Range('E28').Value = '  -8.69%  '

# Row 29: BitcoinCash
$ws.Range('D29').Value = '''126.20'
$ws.Range('E29').Value = '  -6.22%  '

# Row 30: WrappedliquidstakedEther2.0
$ws.Range('D30').Value = '1.826.82'
$ws.Range('E30').Value = '  -7.10%  '

# Row 31: ImmutableX
$ws.Range('D31').Value = '''1.054'
$ws.Range('E31').Value = '  -15.58%  '

# Row 32: HuobiToken
$ws.Range('D32').Value = '''4.088'
$ws.Range('E32').Value = '  -5.24%  '

# Row 33: Filecoin
$ws.Range('D33').Value = '''5.713'
$ws.Range('E33').Value = '  -23.64%  '

# Row 34: Stellar
$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').Value = '''0.08431'
$ws.Range('E34').Value = '  -3.61%  '

# Row 35: WEMIXTOKEN
$ws.Range('B35').Value = 'WEMIXTOKEN'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = '''1.699'
$ws.Range('E35').Value = '  -6.77%  '

# Row 36: Aptos
$ws.Range('D36').Value = '''12.33'
$ws.Range('E36').Value = '  -13.21%  '

# Row 37: InternetComputer(DFINITY)
$ws.Range('D37').Value = '''5.128'
$ws.Range('E37').Value = '  -10.82%  '

# Row 38: Hedera
$ws.Range('D38').Value = '''0.06101'
$ws.Range('E38').Value = '  -10.62%  '

# Row 39: VeChain
$ws.Range('D39').Value = '''0.02228'
$ws.Range('E39').Value = '  -11.29%  '

# Row 40: TrustWalletToken
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '''1.209'
$ws.Range('E40').Value = '  -7.15%  '

# Row 41: Algorand
$ws.Range('D41').Value = '''0.2050'
$ws.Range('E41').Value = '  -8.69%  '

# Row 42: FraxShare
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '''8.188'
$ws.Range('E42').Value = '  -13.30%  '

# Row 43: Frax
$ws.Range('D43').Value = '''1.008'
$ws.Range('E43').Value = '  +1.33%  '

# Row 44: TheSandbox
$ws.Range('D44').Value = '''0.5908'
$ws.Range('E44').Value = '  -10.40%  '

# Row 45: PancakeSwap
$ws.Range('D45').Value = '''3.748'
$ws.Range('E45').Value = '  -4.48%  '

# Row 46: EnergySwap
$ws.Range('D46').Value = '''12.74'
$ws.Range('E46').Value = '  -11.39%  '

# Row 47: Decentraland
$ws.Range('D47').Value = '''0.5621'
$ws.Range('E47').Value = '  -11.93%  '

# Row 48: Quant
$ws.Range('D48').Value = '''121.47'
$ws.Range('E48').Value = '  -8.83%  '

# Row 49: NEARProtocol
$ws.Range('D49').Value = '''1.943'
$ws.Range('E49').Value = '  -10.84%  '

# Row 50: Cronos
$ws.Range('D50').Value = '''0.07017'
$ws.Range('E50').Value = '  -6.74%  '

# Row 51: Aave
$ws.Range('D51').Value = '''73.95'
$ws.Range('E51').Value = '  -8.92%  '
